$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "Not worse"
$ws.Range("C3").Value = "Not worse"
$ws.Range("C4").Value = "A little worse"
$ws.Range("C5").Value = "A little worse"
$ws.Range("C6").Value = "A little worse"
$ws.Range("C7").Value = "A little worse"
$ws.Range("C8").Value = "A little worse"
$ws.Range("C9").Value = "A little worse"
$ws.Range("C10").Value = "A little worse"
$ws.Range("C11").Value = "Not worse"
$ws.Range("C12").Value = "Not worse"
$ws.Range("C13").Value = "Not worse"
$ws.Range("C14").Value = "Not worse"
$ws.Range("C15").Value = "A little worse"

$ws.Range("C16").Select()
